# epv-dependo-refseqs-side-data.xlsx edit
# Commit message: "Refactor (dependo) + content update"
#
# 1) Update existing row 51 (dependo.87-megaderma): F "genome"->"nk", N/O "NULL"->"NK"
# 2) Insert a brand-new row at position 52 for "dependo.88-Megaderma" (a second
#    Megaderma variant), duplicating formatting from (the just-fixed) row 51,
#    and populate it with its own values.
# 3) Update the saved sheet view (pane/selection) to the state captured in the
#    authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: fix up row 51 content/formatting
# ---------------------------------------------------------------------------

# F51: genome -> nk  (style changes from the "genome" look to the "nk" look;
# copy formats from a cell that already carries the "nk" style, e.g. F8)
$ws.Range("F8").Copy()
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("F51").Value2 = "nk"

# N51/O51: NULL -> NK (copy formats from a cell that already has the "NK" look,
# e.g. J51 which is already an "NK" cell in the same row)
$ws.Range("J51").Copy()
$ws.Range("N51:O51").PasteSpecial(-4122)
$ws.Range("N51").Value2 = "NK"
$ws.Range("O51").Value2 = "NK"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: insert a new row at 52 (shifts old rows 52-86 down to 53-87) and
# fill it in with the new "dependo.88-Megaderma" record.
# ---------------------------------------------------------------------------

$ws.Rows(52).Insert()

$ws.Range("A52").Value2 = "dependo.88-Megaderma"
$ws.Range("B52").Value2 = "yes"
$ws.Range("C52").Value2 = 88
$ws.Range("D52").Value2 = "species"
$ws.Range("E52").Value2 = "Megaderma lytra"
$ws.Range("F52").Value2 = "nk"
$ws.Range("G52").Value2 = 1
$ws.Range("H52").Value2 = 1
$ws.Range("I52").Value2 = "no"
$ws.Range("J52").Value2 = "NK"
$ws.Range("K52").Value2 = "NK"
$ws.Range("L52").Value2 = "NK"
$ws.Range("M52").Value2 = "NK"
$ws.Range("N52").Value2 = "NK"
$ws.Range("O52").Value2 = "NK"
$ws.Range("P52").ClearContents()
$ws.Range("Q52").ClearContents()
$ws.Range("R52").Value2 = "ADAMTS2"
$ws.Range("S52").Value2 = "U6"
$ws.Range("T52").Value2 = "NK"
$ws.Range("U52").Value2 = "NK"
$ws.Range("V52").Value2 = "NA"
$ws.Range("W52").Value2 = "dependo.88-Megaderma"
$ws.Range("X52").Value2 = "dependo.88-Megaderma"
$ws.Range("Y52").Value2 = "Parvovirinae"
$ws.Range("Z52").Value2 = "fasta-refseqs-dependo-epv"
$ws.Range("AA52").Value2 = "Dependoparvovirus"

# ---------------------------------------------------------------------------
# Step 3: update the recorded sheet view/selection state (re-affirm the
# frozen header row/column, then leave the final selection on L51 as in the
# authored workbook)
# ---------------------------------------------------------------------------

$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("L51").Select()
